$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts old D:K -> F:M)
$ws.Columns("D:E").Insert()

# Copy cell formatting (number formats/styles) from column F (the old column D, now shifted)
# into the two newly inserted, blank columns D:E so they match the workbook's existing look
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New quarter header dates (30-Nov-2018 and 31-Aug-2018) for the three "Period Ending" rows
$ws.Range("D7").Value = 43434
$ws.Range("E7").Value = 43343
$ws.Range("D38").Value = 43434
$ws.Range("E38").Value = 43343
$ws.Range("D80").Value = 43434
$ws.Range("E80").Value = 43343

# Rows that stay blank in the two new columns (section spacer rows)
$blankRows = 11,16,19,39,40,55,56,67,82,90,95
foreach ($r in $blankRows) {
  $ws.Range("D$r").ClearContents()
  $ws.Range("E$r").ClearContents()
}

# Rows where the two new columns are "NA" (not applicable)
$naRows = 12,29
foreach ($r in $naRows) {
  $ws.Range("D$r").Value = "NA"
  $ws.Range("E$r").Value = "NA"
}

# Numeric values for the two new quarter columns
$newData = @{
  8 = @(9374000, 9948000)
  9 = @(5269000, 5551000)
  10 = @(4105000, 4397000)
  13 = @(0, 0)
  14 = @(0, 0)
  15 = @(0, 0)
  17 = @(8411000, 8614000)
  18 = @(963000, 1334000)
  20 = @(34000, -64000)
  21 = @(1175000, 1451000)
  22 = @(0, 0)
  23 = @(997000, 1270000)
  24 = @(150000, 178000)
  25 = @(0, 0)
  26 = @(847000, 1092000)
  27 = @(847000, 1092000)
  28 = @(0, 0)
  30 = @(0, 0)
  31 = @(0, 0)
  32 = @(-34000, 64000)
  33 = @(847000, 1092000)
  34 = @(0, 0)
  35 = @(847000, 1092000)
  41 = @(3423000, 3282000)
  42 = @(618000, 987000)
  43 = @(4346000, 4330000)
  44 = @(5388000, 5227000)
  45 = @(1791000, 1675000)
  46 = @(15566000, 15501000)
  47 = @(0, 0)
  48 = @(4588000, 4487000)
  49 = @(438000, 438000)
  50 = @(0, 0)
  51 = @(0, 0)
  52 = @(2085000, 2057000)
  53 = @(0, 0)
  54 = @(22677000, 22483000)
  57 = @(2574000, 2333000)
  58 = @(15000, 19000)
  59 = @(4689000, 4356000)
  60 = @(7278000, 6708000)
  61 = @(3466000, 3467000)
  62 = @(3204000, 3316000)
  63 = @(0, 0)
  64 = @(0, 0)
  65 = @(0, 0)
  66 = @(13948000, 13491000)
  68 = @(0, 0)
  69 = @(0, 0)
  70 = @(0, 0)
  71 = @(0, 0)
  72 = @(1810000, 2494000)
  73 = @(0, 0)
  74 = @(0, 0)
  75 = @(0, 0)
  76 = @(8729000, 8992000)
  77 = @(0, 0)
  81 = @(847000, 1092000)
  83 = @(178000, 181000)
  84 = @(0, 0)
  85 = @(0, 0)
  86 = @(0, 0)
  87 = @(0, 0)
  88 = @(0, 0)
  89 = @(1524000, 1301000)
  91 = @(-287000, -343000)
  92 = @(0, 0)
  93 = @(0, 0)
  94 = @(88000, -333000)
  96 = @(-318000, -320000)
  97 = @(0, 0)
  98 = @(0, 0)
  99 = @(0, 0)
  100 = @(-1458000, -1832000)
  101 = @(-13000, -103000)
  102 = @(141000, -967000)
}
foreach ($r in $newData.Keys) {
  $vals = $newData[$r]
  $ws.Range("D$r").Value = $vals[0]
  $ws.Range("E$r").Value = $vals[1]
}

# Data correction: row 91 column J (Changes In Other Operating Activities, FY2017 Q3)
# was restated from -658000 to -329000
$ws.Range("J91").Value = -329000

